$d = $word.ActiveDocument
$p = $d.Paragraphs(4)
$rng = $p.Range
$xml = '<w:p w14:paraId="01A55588" w14:textId="09B74061" w:rsidR="00361123" w:rsidRPr="00FC33D9" w:rsidRDefault="00361123"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:t xml:space="preserve">7.Чтобы создать внутреннюю ссылку необходимо выбрать </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>слой</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> с которым мы будем работать при перемещении</w:t></w:r><w:r w:rsidRPr="00361123"><w:t xml:space="preserve">, </w:t></w:r><w:r><w:t>перейти во вкладку прототип</w:t></w:r><w:r w:rsidRPr="00361123"><w:t>,</w:t></w:r><w:r><w:t xml:space="preserve"> во вкладке </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Interactions</w:t></w:r><w:r w:rsidRPr="00361123"><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>выбрать событие при котором будет происходить перемещение между фреймами и выбрать фрейм</w:t></w:r><w:r w:rsidRPr="00361123"><w:t>,</w:t></w:r><w:r><w:t xml:space="preserve"> также можно просто нажать прототип</w:t></w:r><w:r w:rsidRPr="00361123"><w:t>,</w:t></w:r><w:r><w:t xml:space="preserve"> навести курсор ближе к концу слоя и </w:t></w:r><w:r w:rsidR="00FC33D9"><w:t>потянуть плюсик на фрейм к которому нам нужно переместиться</w:t></w:r><w:r w:rsidR="00FC33D9" w:rsidRPr="00FC33D9"><w:t>,</w:t></w:r><w:r w:rsidR="00FC33D9"><w:t xml:space="preserve"> там в дальнейшем можно выбирать уже эффекты переходов</w:t></w:r><w:r w:rsidR="00FC33D9" w:rsidRPr="00FC33D9"><w:t>,</w:t></w:r><w:r w:rsidR="00FC33D9"><w:t xml:space="preserve"> скорость и т.д.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>https://www.figma.com/file/eosiHLHRYQpnqyEwWlMARQ/Untitled?type=design&amp;t=N5iUg1F5tFdtISVi-6</w:t></w:r></w:p>'
$rng.InsertXML($xml)
Write-Output ("Paragraphs after: " + $d.Paragraphs.Count)
